$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# J2 changes from "002" to "001" (kept as text, matching K2's existing text "001")
$ws.Range("K2").Copy()
$ws.Range("J2").PasteSpecial(-4163)

$ws.Range("N2").Value = "2017-12-31 00:00:00"
$ws.Range("O2").Value = 32260076.08
$ws.Range("P2").Value = 381909527.23
$ws.Range("Q2").Value = 347682544.98
$ws.Range("R2").Value = 5.3248569944
$ws.Range("S2").Value = 254156707.02
$ws.Range("T2").Value = 254156707.02
$ws.Range("U2").Value = 2.5852564758
$ws.Range("V2").Value = 41230517.34
$ws.Range("W2").Value = 26339237.82
$ws.Range("X2").Value = 1732933.18
$ws.Range("Y2").Value = 36803822.04
$ws.Range("Z2").Value = 36623900.99
$ws.Range("AA2").Value = 4401688.82
$ws.Range("AG2").Value = 452609.16
$ws.Range("AP2").Value = 7.0611582777
$ws.Range("AQ2").Value = 45.450067394206
$ws.Range("AR2").Value = 45.427427249849
$ws.Range("AS2").Value = 31064058.97
$ws.Range("AT2").Value = 41.818582310869
